# "added ssi to dic"
#
# 1) Insert a new list item "SSI; surgical site infections" (same level /
#    numbering / run-formatting as the neighbouring "RTI ..." / "UTI ..."
#    items) right after the "UTI; urinary tract infection" paragraph, and
#    carry the Word "last edit" bookmark (_GoBack) onto it.
# 2) The _GoBack bookmark used to sit in the "Additionally unique_id_4both
#    was added in ..." paragraph, splitting one sentence into two runs.
#    Since the bookmark moved (step 1), that split collapses back into a
#    single run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: add the "SSI; surgical site infections" bullet after "UTI; ..."
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*UTI; urinary tract infection*") {

        # Clone the UTI paragraph's formatting/mark by inserting a new
        # (empty) paragraph straight after it ...
        $p.Range.InsertParagraphAfter() | Out-Null
        $newPara = $p.Next()

        # ... then stamp the new paragraph with the exact OOXML Word would
        # produce for this bullet (keeps the run free of an explicit rPr,
        # matching its sibling items, and attaches the _GoBack bookmark
        # that Word leaves at the most recent edit location).
        $paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="McgjyyAdvTTb5929f4c" w:hAnsi="McgjyyAdvTTb5929f4c" w:cs="McgjyyAdvTTb5929f4c"/><w:color w:val="131413"/></w:rPr></w:pPr><w:r><w:t>SSI; surgical site infections</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
        $newPara.Range.InsertXML($paraXml) | Out-Null

        break
    }
}

# ---------------------------------------------------------------------
# Step 2: drop the old _GoBack bookmark and rejoin the sentence it split
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*checked by hand, due to removals*") {

        $fixedXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3777215A" w14:textId="2873DDFC" w:rsidR="000744B3" w:rsidRDefault="000744B3" w:rsidP="000744B3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:contextualSpacing w:val="0"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>Additionally</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> unique_id</w:t></w:r><w:r w:rsidR="008311FF"><w:t>_4both</w:t></w:r><w:r><w:t xml:space="preserve"> was </w:t></w:r><w:r w:rsidR="0064515A"><w:t xml:space="preserve">added in by hand throughout the extraction process, and subsequently </w:t></w:r><w:r><w:t xml:space="preserve">checked by hand, due to removals of rows in the data checking and cleaning process there may be gaps in numbers. For </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>example</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>…</w:t></w:r></w:p>'
        $p.Range.InsertXML($fixedXml) | Out-Null

        break
    }
}
